# "Update automatico via Actualizar 06-08-2020 03-06-43"
# Appends two new daily rows (6/6/2020 and 7/6/2020) to the
# "Condicion_Pacientes" table on Hoja1, and grows the table/autofilter
# range from A1:I70 to A1:I72 to include them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lo = $ws.ListObjects.Item("Condicion_Pacientes")

# Grow the table (and its autofilter) down to row 72 first, so the new
# rows become part of the table before we populate them.
$lo.Resize($ws.Range("A1:I72"))

# Copy the formatting (styles) of the last existing data row down onto
# the two new rows, matching how the table's row style is carried
# forward for newly appended rows.
$ws.Range("A70:F70").Copy() | Out-Null
$ws.Range("A71:F71").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A70:F70").Copy() | Out-Null
$ws.Range("A72:F72").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# Row 71: 6/6/2020
$ws.Range("A71").Formula = "=+Condicion_Pacientes[[#This Row],[día]]&""/""&Condicion_Pacientes[[#This Row],[mes]]&""/""&Condicion_Pacientes[[#This Row],[año]]"
$ws.Range("B71").Value = 6
$ws.Range("C71").Value = 6
$ws.Range("D71").Value = 2020
$ws.Range("E71").Value = 1630
$ws.Range("F71").Value = 307

# Row 72: 7/6/2020
$ws.Range("A72").Formula = "=+Condicion_Pacientes[[#This Row],[día]]&""/""&Condicion_Pacientes[[#This Row],[mes]]&""/""&Condicion_Pacientes[[#This Row],[año]]"
$ws.Range("B72").Value = 7
$ws.Range("C72").Value = 6
$ws.Range("D72").Value = 2020
$ws.Range("E72").Value = 1749
$ws.Range("F72").Value = 263

# Move the on-screen selection to where the user last left off entering
# data (matches the saved view state).
$ws.Range("F73").Select() | Out-Null
